$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# ---------------------------------------------------------------------
# Add more entry points (new survey rows) below the existing data table.
# Rows 41-44 are brand new rows; rows 45-46 were already present as blank
# placeholder rows (pre-formatted with the shaded style) and now get
# filled in with data as well.
# ---------------------------------------------------------------------

$ws.Range("A41").Value = "Male"
$ws.Range("B41").Value = 16
$ws.Range("C41").Value = 118
$ws.Range("D41").Value = "told"

$ws.Range("A42").Value = "Female"
$ws.Range("B42").Value = 19
$ws.Range("C42").Value = 43
$ws.Range("D42").Value = "seen"

$ws.Range("A43").Value = "Male"
$ws.Range("B43").Value = 24
$ws.Range("C43").Value = 81
$ws.Range("D43").Value = "seen"

$ws.Range("A44").Value = "Male"
$ws.Range("B44").Value = 19
$ws.Range("C44").Value = 62
$ws.Range("D44").Value = "seen"

$ws.Range("A45").Value = "Male"
$ws.Range("B45").Value = 19
$ws.Range("C45").Value = 112
$ws.Range("D45").Value = "seen"

$ws.Range("A46").Value = "Female"
$ws.Range("B46").Value = 19
$ws.Range("C46").Value = 107
$ws.Range("D46").Value = "seen"

# ---------------------------------------------------------------------
# Reflect the newly-added rows in the window/view state: scroll the
# sheet down near the bottom of the data and select the last entered
# row's last cell, same as the author would have after typing the data.
# ---------------------------------------------------------------------

$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("D40").Select()
